$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'28.086.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  +1.94%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'1.910.06"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  +1.81%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.89%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'316.76"
$ws.Range("D5").Style = "Normal"
$ws.Range("E6").Value = "'  -0.80%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'0.4825"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  +0.58%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.3824"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  +0.97%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.07367"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -0.29%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.9340"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -0.77%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'20.78"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  +0.24%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.07835"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -0.49%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'1.913.83"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  +1.57%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'5.503"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  +0.83%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'6.617"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  +0.17%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'91.30"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'  +0.04%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'1.009"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -0.83%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'0.000008839"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -1.57%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'1.007"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -0.73%  "
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'28.112.20"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'  +1.92%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'14.82"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -0.94%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'5.155"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  +0.15%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'2.154.13"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  +1.43%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'10.88"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  +0.88%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'156.71"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  +1.90%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'1.920"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -2.39%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'18.56"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -0.15%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'2.105"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  +4.00%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("D29").Value = "'116.41"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  +0.24%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'4.962"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -1.18%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'0.08903"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -0.41%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'3.363"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  +1.02%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'1.244"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").Value = "'0.7675"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  +2.14%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("D35").Value = "'4.681"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  +1.25%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "'2.614"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -3.15%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.02044"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "'  -1.53%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("D39").Value = "'0.05301"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -0.26%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.5498"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  +2.33%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'2.995"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -0.52%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'7.016"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").Value = "'0.1524"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -0.10%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'8.452"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +0.06%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'10.74"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  +0.35%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.4837"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -0.31%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'107.13"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  +3.89%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.007"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -0.84%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'1.656"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -0.49%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'68.39"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  +1.48%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.06096"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "'  -0.19%  "
$ws.Range("E51").Style = "Normal"
